# Auto-generated edit script applying numeric updates to the Leve profit sheets
# (currentAveragePrice / LevePrice / LeveProfit columns) per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 42860.105
$ws.Range("I98").Value = 1281.5
$ws.Range("J98").Value = 89058.55499999999
$ws.Range("K98").Value = 1281.5
$ws.Range("L98").Value = 89058.55499999999
$ws.Range("M98").Value = 216.5
$ws.Range("N98").Value = -92054.55499999999

$ws.Range("H114").Value = 39372.668
$ws.Range("J114").Value = 39372.668
$ws.Range("L114").Value = 39372.668
$ws.Range("N114").Value = -48050.668

$ws.Range("H120").Value = 49716.668
$ws.Range("J120").Value = 49716.668
$ws.Range("L120").Value = 49716.668
$ws.Range("N120").Value = -59392.668

$ws.Range("H122").Value = 42860.105
$ws.Range("I122").Value = 1281.5
$ws.Range("J122").Value = 89058.55499999999
$ws.Range("K122").Value = 3844.5
$ws.Range("L122").Value = 267175.665
$ws.Range("M122").Value = -1394.5
$ws.Range("N122").Value = -272075.665

$ws.Range("H128").Value = 36793.668
$ws.Range("J128").Value = 36793.668
$ws.Range("L128").Value = 36793.668
$ws.Range("N128").Value = -46753.668

$ws.Range("H129").Value = 1306.7073
$ws.Range("J129").Value = 1213.3939
$ws.Range("L129").Value = 3640.1817
$ws.Range("N129").Value = -13640.1817

$ws.Range("H137").Value = 3645.8655
$ws.Range("I137").Value = 1025.3334
$ws.Range("J137").Value = 4432.025
$ws.Range("K137").Value = 3076.0002
$ws.Range("L137").Value = 13296.075
$ws.Range("M137").Value = -526.0001999999999
$ws.Range("N137").Value = -18396.075

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 31159.797
$ws.Range("I32").Value = 32594.94
$ws.Range("J32").Value = 23186.777
$ws.Range("K32").Value = 32594.94
$ws.Range("L32").Value = 23186.777
$ws.Range("M32").Value = -32307.94
$ws.Range("N32").Value = -23760.777

$ws.Range("H80").Value = 59102
$ws.Range("J80").Value = 59102
$ws.Range("L80").Value = 59102
$ws.Range("N80").Value = -61098

$ws.Range("H83").Value = 59102
$ws.Range("J83").Value = 59102
$ws.Range("L83").Value = 177306
$ws.Range("N83").Value = -187290

$ws.Range("H125").Value = 37925.168
$ws.Range("J125").Value = 37925.168
$ws.Range("L125").Value = 37925.168
$ws.Range("N125").Value = -47765.168

$ws.Range("H132").Value = 1591.5238
$ws.Range("I132").Value = 1128.5
$ws.Range("J132").Value = 2517.5715
$ws.Range("K132").Value = 3385.5
$ws.Range("L132").Value = 7552.7145
$ws.Range("M132").Value = -855.5
$ws.Range("N132").Value = -12612.7145

$ws.Range("H133").Value = 28595
$ws.Range("J133").Value = 28595
$ws.Range("L133").Value = 28595
$ws.Range("N133").Value = -33655

$ws.Range("H134").Value = 51178
$ws.Range("J134").Value = 51178
$ws.Range("L134").Value = 51178
$ws.Range("N134").Value = -61318

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2290.261
$ws.Range("I107").Value = 2156.9167
$ws.Range("J107").Value = 2435.7273
$ws.Range("K107").Value = 2156.9167
$ws.Range("L107").Value = 2435.7273
$ws.Range("M107").Value = -236.9167000000002
$ws.Range("N107").Value = -6275.7273

$ws.Range("H119").Value = 47753
$ws.Range("J119").Value = 47753
$ws.Range("L119").Value = 47753
$ws.Range("N119").Value = -57429

$ws.Range("H123").Value = 25000
$ws.Range("J123").Value = 25000
$ws.Range("L123").Value = 25000
$ws.Range("N123").Value = -34800

$ws.Range("H124").Value = 47992
$ws.Range("J124").Value = 47992
$ws.Range("L124").Value = 47992
$ws.Range("N124").Value = -57812

$ws.Range("H132").Value = 52460
$ws.Range("J132").Value = 52460
$ws.Range("L132").Value = 52460
$ws.Range("N132").Value = -62580

$ws.Range("H137").Value = 44866.332
$ws.Range("J137").Value = 44866.332
$ws.Range("L137").Value = 44866.332
$ws.Range("N137").Value = -55066.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 47461.2
$ws.Range("J20").Value = 47461.2
$ws.Range("L20").Value = 47461.2
$ws.Range("N20").Value = -47933.2

$ws.Range("H30").Value = 47461.2
$ws.Range("J30").Value = 47461.2
$ws.Range("L30").Value = 47461.2
$ws.Range("N30").Value = -47643.2

$ws.Range("H31").Value = 199264.97
$ws.Range("I31").Value = 2548.8333
$ws.Range("J31").Value = 238608.2
$ws.Range("K31").Value = 2548.8333
$ws.Range("L31").Value = 238608.2
$ws.Range("M31").Value = -2253.8333
$ws.Range("N31").Value = -239198.2

$ws.Range("H34").Value = 199264.97
$ws.Range("I34").Value = 2548.8333
$ws.Range("J34").Value = 238608.2
$ws.Range("K34").Value = 2548.8333
$ws.Range("L34").Value = 238608.2
$ws.Range("M34").Value = -2346.8333
$ws.Range("N34").Value = -239012.2

$ws.Range("H53").Value = 10684
$ws.Range("J53").Value = 10684
$ws.Range("L53").Value = 10684
$ws.Range("N53").Value = -11898

$ws.Range("H108").Value = 20000
$ws.Range("J108").Value = 20000
$ws.Range("L108").Value = 20000
$ws.Range("N108").Value = -27680

$ws.Range("H110").Value = 46463.668
$ws.Range("J110").Value = 46463.668
$ws.Range("L110").Value = 46463.668
$ws.Range("N110").Value = -54643.668

$ws.Range("H111").Value = 48698
$ws.Range("J111").Value = 48698
$ws.Range("L111").Value = 48698
$ws.Range("N111").Value = -56878

$ws.Range("H112").Value = 28025.143
$ws.Range("J112").Value = 28025.143
$ws.Range("L112").Value = 28025.143
$ws.Range("N112").Value = -30979.143

$ws.Range("H128").Value = 47461.2
$ws.Range("J128").Value = 47461.2
$ws.Range("L128").Value = 47461.2
$ws.Range("N128").Value = -57421.2

$ws.Range("H132").Value = 103090.29
$ws.Range("I132").Value = 2629.889
$ws.Range("K132").Value = 7889.667
$ws.Range("M132").Value = -5359.667

$ws.Range("H138").Value = 47866.332
$ws.Range("J138").Value = 47866.332
$ws.Range("L138").Value = 47866.332
$ws.Range("N138").Value = -58146.332

$ws.Range("H141").Value = 18223.25
$ws.Range("J141").Value = 18223.25
$ws.Range("L141").Value = 18223.25
$ws.Range("N141").Value = -28583.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1771
$ws.Range("I102").Value = 1593.3846
$ws.Range("J102").Value = 2232.8
$ws.Range("K102").Value = 1593.3846
$ws.Range("L102").Value = 2232.8
$ws.Range("M102").Value = 28.61539999999991
$ws.Range("N102").Value = -5476.8

$ws.Range("H110").Value = 33503
$ws.Range("J110").Value = 33503
$ws.Range("L110").Value = 33503
$ws.Range("N110").Value = -41683

$ws.Range("H122").Value = 2028.5714
$ws.Range("I122").Value = 2275
$ws.Range("J122").Value = 1700
$ws.Range("K122").Value = 6825
$ws.Range("L122").Value = 5100
$ws.Range("M122").Value = -4375
$ws.Range("N122").Value = -10000

$ws.Range("H130").Value = 46474.5
$ws.Range("J130").Value = 46474.5
$ws.Range("L130").Value = 46474.5
$ws.Range("N130").Value = -56514.5

$ws.Range("H132").Value = 2918.761
$ws.Range("I132").Value = 1205.1818
$ws.Range("J132").Value = 7268.615
$ws.Range("K132").Value = 3615.5454
$ws.Range("L132").Value = 21805.845
$ws.Range("M132").Value = -1085.5454
$ws.Range("N132").Value = -26865.845

$ws.Range("H135").Value = 39158.2
$ws.Range("J135").Value = 39158.2
$ws.Range("L135").Value = 39158.2
$ws.Range("N135").Value = -49298.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H111").Value = 35549.168
$ws.Range("J111").Value = 35549.168
$ws.Range("L111").Value = 35549.168
$ws.Range("N111").Value = -43729.168

$ws.Range("H112").Value = 26654.666
$ws.Range("J112").Value = 28985.6
$ws.Range("L112").Value = 28985.6
$ws.Range("N112").Value = -31939.6

$ws.Range("H121").Value = 25305.666
$ws.Range("J121").Value = 25305.666
$ws.Range("L121").Value = 25305.666
$ws.Range("N121").Value = -28799.666

$ws.Range("H132").Value = 3298.4888
$ws.Range("I132").Value = 1984.9259
$ws.Range("J132").Value = 5268.8335
$ws.Range("K132").Value = 5954.7777
$ws.Range("L132").Value = 15806.5005
$ws.Range("M132").Value = -3424.7777
$ws.Range("N132").Value = -20866.5005

$ws.Range("H136").Value = 2146.641
$ws.Range("I136").Value = 1621.5518
$ws.Range("J136").Value = 3669.4
$ws.Range("K136").Value = 4864.6554
$ws.Range("L136").Value = 11008.2
$ws.Range("M136").Value = -2314.6554
$ws.Range("N136").Value = -16108.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 855.5
$ws.Range("I107").Value = 807.3333
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 2421.9999
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = -501.9998999999998
$ws.Range("N107").Value = -6840

$ws.Range("H110").Value = 27634.8
$ws.Range("J110").Value = 27634.8
$ws.Range("L110").Value = 27634.8
$ws.Range("N110").Value = -35814.8

$ws.Range("H112").Value = 29377.75
$ws.Range("J112").Value = 29377.75
$ws.Range("L112").Value = 29377.75
$ws.Range("N112").Value = -32331.75

$ws.Range("H118").Value = 30299.4
$ws.Range("J118").Value = 34124.25
$ws.Range("L118").Value = 34124.25
$ws.Range("N118").Value = -37438.25

$ws.Range("H137").Value = 26508.334
$ws.Range("J137").Value = 26508.334
$ws.Range("L137").Value = 26508.334
$ws.Range("N137").Value = -36708.334

Write-Host "Applied profit updates to ALC, ARM, BSM, CRP, GSM, LTW, WVR sheets"
